# Generate Report for Handoff
#
# Updates the localization-status report for file
# "4a452103-27ad-4a63-a2bc-f2f24f1b01bc.md" in both the Overview sheet
# and the two per-locale sheets (zh-cn, de-de): the item moved from
# "Handed back: in sync with en-US" to "Ready for handoff", the handoff
# timestamps were refreshed, and an error detail was recorded because the
# handback file is behind the latest source.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4032585cd81735975657acdcddfb906b8c926edc/e2e/4a452103-27ad-4a63-a2bc-f2f24f1b01bc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41dc2ba87688ad1c1efb3f819e7ef78ff424fec2/e2e/4a452103-27ad-4a63-a2bc-f2f24f1b01bc.md."

# --- Overview sheet: row 3 is the 4a452103-27ad-4a63-a2bc-f2f24f1b01bc.md entry ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = "2016-08-30 15:03:30"

# --- zh-cn sheet: row 3 is the 4a452103-27ad-4a63-a2bc-f2f24f1b01bc.md entry ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = "2016-08-30 15:03:25"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: row 3 is the 4a452103-27ad-4a63-a2bc-f2f24f1b01bc.md entry ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = "2016-08-30 15:03:30"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
